# Insert a new data row at row 54 (pushes existing rows 54..137 down to 55..138)
# and populate it with the new record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(54).Insert()

$ws.Cells.Item(54, 1).Value  = 11
$ws.Cells.Item(54, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(54, 3).Value  = "Bíobío"
$ws.Cells.Item(54, 4).Value  = 44967
$ws.Cells.Item(54, 5).Value  = 8
$ws.Cells.Item(54, 6).Value  = "Fruta"
$ws.Cells.Item(54, 7).Value  = 100101
$ws.Cells.Item(54, 8).Value  = "Berries"
$ws.Cells.Item(54, 9).Value  = 100101001
$ws.Cells.Item(54, 10).Value = "Arándano (blue)"
$ws.Cells.Item(54, 11).Value = "Sin especificar"
$ws.Cells.Item(54, 12).Value = "Primera"
$ws.Cells.Item(54, 13).Value = 150
$ws.Cells.Item(54, 14).Value = 3500
$ws.Cells.Item(54, 15).Value = 3600
$ws.Cells.Item(54, 16).Value = 3533
$ws.Cells.Item(54, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(54, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(54, 19).Value = 1766
$ws.Cells.Item(54, 20).Value = 2
